# "fixed error in experiment format"
#
# The "Survey 2" sheet was missing the "Pseudo-Random Question Width:"
# label row that "Survey 1" already has right after the "Questions Per
# Page:" row. This inserts the missing row (shifting everything below it
# down by one, matching Survey 1's layout) and re-activates "Survey 2"
# as the selected/visible sheet.

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Survey 2")

# Insert the missing label row right above the old row 8 ("Width:" row),
# pushing the rest of the sheet down by one row - this mirrors the
# structure already present on "Survey 1".
$ws2.Rows.Item(8).Insert()
$ws2.Range("A8").Value = "Pseudo-Random Question Width:"

# The fix was made while reviewing "Survey 2", so it ends up the active
# (selected) sheet/tab when the workbook is saved.
$ws2.Activate()
